{"js": "// Turn the two placeholder label paragraphs into Mustache-style template\n// fields:\n//   \"Document_id = \"          -> \"Name = {{name}}\"\n//   \"Document_template_id = \" -> \"Position = {{position}}\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n  if (text === \"Document_id = \") {\n    paragraph.insertText(\"Name = {{name}}\", \"Replace\");\n  } else if (text === \"Document_template_id = \") {\n    paragraph.insertText(\"Position = {{position}}\", \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the two placeholder label paragraphs with their templated equivalents.\n# \"Document_id = \"          -> \"Name = {{name}}\"\n# \"Document_template_id = \" -> \"Position = {{position}}\"\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $t = $r.Text.TrimEnd(\"`r\")\n    if ($t -eq \"Document_id = \") {\n        $r.Text = \"Name = {{name}}\"\n    }\n    elseif ($t -eq \"Document_template_id = \") {\n        $r.Text = \"Position = {{position}}\"\n    }\n}\n"}
